$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:F (runs, balls, fours, sixes) for rows 2-10
$data = @{
    2  = @("5", "6", "0", "0")
    3  = @("0", "1", "0", "0")
    4  = @("5", "2", "1", "0")
    5  = @("21", "18", "2", "0")
    6  = @("10", "6", "1", "0")
    7  = @("30", "27", "2", "1")
    9  = @("10", "10", "0", "1")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt 4; $i++) {
        $cell = $ws.Cells.Item($row, 3 + $i)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
    }
}
